$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(32, 1).Value = "2025-04-28 23:29:29"
$ws.Cells.Item(32, 2).Value = 192
